$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark rows as "Processed" in column C: the 5 rows (11-15) that previously
# had no value in column C, plus every 10th "Capricious" word-group row
# (37, 47, 57, ... 1027) whose column B reads "Well-meaning and kindly."
$rows = @(11, 12, 13, 14, 15, 37, 47, 57, 67, 77, 87, 97, 107, 117, 127, 137, 147, 157, 167, 177, 187, 197, 207, 217, 227, 237, 247, 257, 267, 277, 287, 297, 307, 317, 327, 337, 347, 357, 367, 377, 387, 397, 407, 417, 427, 437, 447, 457, 467, 477, 487, 497, 507, 517, 527, 537, 547, 557, 567, 577, 587, 597, 607, 617, 627, 637, 647, 657, 667, 677, 687, 697, 707, 717, 727, 737, 747, 757, 767, 777, 787, 797, 807, 817, 827, 837, 847, 857, 867, 877, 887, 897, 907, 917, 927, 937, 947, 957, 967, 977, 987, 997, 1007, 1017, 1027)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = "Processed"
}
